$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 1.865848376575912
$ws.Range("D2").Value2 = 0.1779072435960671
$ws.Range("E2").Value2 = 0.8085047307785374
$ws.Range("F2").Value2 = 3.422755640162137
$ws.Range("G2").Value2 = 0.002447628687513673
$ws.Range("L2").Value2 = 0.81099605637894
$ws.Range("N2").Value2 = 1.535965110662168
$ws.Range("B3").Value2 = 1.773780148288211
$ws.Range("D3").Value2 = 0.1665998921905327
$ws.Range("E3").Value2 = 0.7057250909312813
$ws.Range("F3").Value2 = 3.196875946888753
$ws.Range("G3").Value2 = 0.002459791236997574
$ws.Range("L3").Value2 = 0.732628589055679
$ws.Range("N3").Value2 = 1.561931636863399
$ws.Range("B4").Value2 = 1.71846754994408
$ws.Range("D4").Value2 = 0.1598496732897416
$ws.Range("E4").Value2 = 0.6426644111768667
$ws.Range("F4").Value2 = 3.061510461391663
$ws.Range("G4").Value2 = 0.002467620432423804
$ws.Range("L4").Value2 = 0.6850079866848091
$ws.Range("N4").Value2 = 1.578622067910064
$ws.Range("B5").Value2 = 1.696228960981784
$ws.Range("D5").Value2 = 0.1571456858269187
$ws.Range("E5").Value2 = 0.6169703906518862
$ws.Range("F5").Value2 = 3.007157015669947
$ws.Range("G5").Value2 = 0.002470902274912622
$ws.Range("L5").Value2 = 0.665722229461096
$ws.Range("N5").Value2 = 1.585611290619365
$ws.Range("B6").Value2 = 1.692554357502559
$ws.Range("D6").Value2 = 0.1566994692101957
$ws.Range("E6").Value2 = 0.6127039143550093
$ws.Range("F6").Value2 = 2.998179753380924
$ws.Range("G6").Value2 = 0.002471452756889152
$ws.Range("L6").Value2 = 0.6625269488687309
$ws.Range("N6").Value2 = 1.586783182736678
$ws.Range("B7").Value2 = 1.718166416579095
$ws.Range("D7").Value2 = 0.1598130190839271
$ws.Range("E7").Value2 = 0.6423178879346096
$ws.Range("F7").Value2 = 3.060774191231928
$ws.Range("G7").Value2 = 0.002467664321812224
$ws.Range("L7").Value2 = 0.6847474127955024
$ws.Range("N7").Value2 = 1.578715566950677
$ws.Range("B8").Value2 = 1.833847774406706
$ws.Range("D8").Value2 = 0.1739675081772845
$ws.Range("E8").Value2 = 0.7730517883892247
$ws.Range("F8").Value2 = 3.344165488145961
$ws.Range("G8").Value2 = 0.002451747653430293
$ws.Range("L8").Value2 = 0.7838686347086252
$ws.Range("N8").Value2 = 1.544763309635243
$ws.Range("B9").Value2 = 2.070580403298095
$ws.Range("D9").Value2 = 0.203328960512124
$ws.Range("E9").Value2 = 1.030168942533379
$ws.Range("F9").Value2 = 3.927564504137735
$ws.Range("G9").Value2 = 0.00242337824403009
$ws.Range("L9").Value2 = 0.9824341925690874
$ws.Range("N9").Value2 = 1.484111890929466
$ws.Range("B10").Value2 = 2.25086606992744
$ws.Range("D10").Value2 = 0.2259901043891546
$ws.Range("E10").Value2 = 1.22008710816462
$ws.Range("F10").Value2 = 4.37492039594548
$ws.Range("G10").Value2 = 0.002404234569912922
$ws.Range("L10").Value2 = 1.1312343501653
$ws.Range("N10").Value2 = 1.443168672389355
$ws.Range("B11").Value2 = 2.33434013894481
$ws.Range("D11").Value2 = 0.2365609245833298
$ws.Range("E11").Value2 = 1.306833463085212
$ws.Range("F11").Value2 = 4.582924738272141
$ws.Range("G11").Value2 = 0.002395887172824298
$ws.Range("L11").Value2 = 1.199642418687347
$ws.Range("N11").Value2 = 1.425328963656437
$ws.Range("B12").Value2 = 2.366165874859973
$ws.Range("D12").Value2 = 0.2406036190827479
$ws.Range("E12").Value2 = 1.339743195633559
$ws.Range("F12").Value2 = 4.662372881434237
$ws.Range("G12").Value2 = 0.00239277756742516
$ws.Range("L12").Value2 = 1.225656767103374
$ws.Range("N12").Value2 = 1.418686659591039
$ws.Range("B13").Value2 = 2.359301935387407
$ws.Range("D13").Value2 = 0.239731150883074
$ws.Range("E13").Value2 = 1.332652626196278
$ws.Range("F13").Value2 = 4.645231436473239
$ws.Range("G13").Value2 = 0.002393445000908825
$ws.Range("L13").Value2 = 1.220049127841548
$ws.Range("N13").Value2 = 1.420112159648067
$ws.Range("B14").Value2 = 2.336954100053276
$ws.Range("D14").Value2 = 0.2368927107144145
$ws.Range("E14").Value2 = 1.309539688703069
$ws.Range("F14").Value2 = 4.589447127269182
$ws.Range("G14").Value2 = 0.002395630317354473
$ws.Range("L14").Value2 = 1.201780398565802
$ws.Range("N14").Value2 = 1.424780227023327
$ws.Range("B15").Value2 = 2.323293700502745
$ws.Range("D15").Value2 = 0.2351593236180065
$ws.Range("E15").Value2 = 1.295390571190438
$ws.Range("F15").Value2 = 4.555367407064637
$ws.Range("G15").Value2 = 0.002396975560100329
$ws.Range("L15").Value2 = 1.19060475532342
$ws.Range("N15").Value2 = 1.427654303442822
$ws.Range("B16").Value2 = 2.245440736799139
$ws.Range("D16").Value2 = 0.2253047258373897
$ws.Range("E16").Value2 = 1.21442596879541
$ws.Range("F16").Value2 = 4.361420279749041
$ws.Range("G16").Value2 = 0.002404787311634073
$ws.Range("L16").Value2 = 1.126778716818876
$ws.Range("N16").Value2 = 1.444350367226843
$ws.Range("B17").Value2 = 2.198058768755061
$ws.Range("D17").Value2 = 0.219327834198765
$ws.Range("E17").Value2 = 1.16485400322631
$ws.Range("F17").Value2 = 4.243616579564332
$ws.Range("G17").Value2 = 0.002409671672851896
$ws.Range("L17").Value2 = 1.087811939150583
$ws.Range("N17").Value2 = 1.454794230513546
$ws.Range("B18").Value2 = 2.170943089784771
$ws.Range("D18").Value2 = 0.2159146319574177
$ws.Range("E18").Value2 = 1.136373982040453
$ws.Range("F18").Value2 = 4.176280433499585
$ws.Range("G18").Value2 = 0.002412515056026803
$ws.Range("L18").Value2 = 1.065466473849256
$ws.Range("N18").Value2 = 1.460875189947885
$ws.Range("B19").Value2 = 2.161785574036116
$ws.Range("D19").Value2 = 0.2147631348946675
$ws.Range("E19").Value2 = 1.126736407241793
$ws.Range("F19").Value2 = 4.153552925435122
$ws.Range("G19").Value2 = 0.002413483637846827
$ws.Range("L19").Value2 = 1.057912031992544
$ws.Range("N19").Value2 = 1.462946788026152
$ws.Range("B20").Value2 = 2.203088421714313
$ws.Range("D20").Value2 = 0.2199615292425108
$ws.Range("E20").Value2 = 1.17012759121917
$ws.Range("F20").Value2 = 4.25611312368585
$ws.Range("G20").Value2 = 0.002409148206602146
$ws.Range("L20").Value2 = 1.09195301896284
$ws.Range("N20").Value2 = 1.453674810735471
$ws.Range("B21").Value2 = 2.343512293697984
$ws.Range("D21").Value2 = 0.2377253337292586
$ws.Range("E21").Value2 = 1.316326786055214
$ws.Range("F21").Value2 = 4.605813558064881
$ws.Range("G21").Value2 = 0.002394987047063857
$ws.Range("L21").Value2 = 1.207143337557909
$ws.Range("N21").Value2 = 1.423406028045783
$ws.Range("B22").Value2 = 2.436548804507026
$ws.Range("D22").Value2 = 0.2495677318840421
$ws.Range("E22").Value2 = 1.412235662299338
$ws.Range("F22").Value2 = 4.838351623413985
$ws.Range("G22").Value2 = 0.002386031085473128
$ws.Range("L22").Value2 = 1.283069433812102
$ws.Range("N22").Value2 = 1.404283609857433
$ws.Range("B23").Value2 = 2.386776073310216
$ws.Range("D23").Value2 = 0.2432252282529532
$ws.Range("E23").Value2 = 1.361010938489699
$ws.Range("F23").Value2 = 4.713865057429587
$ws.Range("G23").Value2 = 0.002390783863872655
$ws.Range("L23").Value2 = 1.242485262935929
$ws.Range("N23").Value2 = 1.4144291304676
$ws.Range("B24").Value2 = 2.200814127181388
$ws.Range("D24").Value2 = 0.2196749642522775
$ws.Range("E24").Value2 = 1.167743341923057
$ws.Range("F24").Value2 = 4.250462221508656
$ws.Range("G24").Value2 = 0.002409384755543507
$ws.Range("L24").Value2 = 1.090080659552257
$ws.Range("N24").Value2 = 1.45418066170841
$ws.Range("B25").Value2 = 2.005446237062586
$ws.Range("D25").Value2 = 0.1952028624876334
$ws.Range("E25").Value2 = 0.9604753759377047
$ws.Range("F25").Value2 = 3.766588364740215
$ws.Range("G25").Value2 = 0.002430752023853131
$ws.Range("L25").Value2 = 0.9282319864749695
$ws.Range("N25").Value2 = 1.499884743834691
